$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers (losing the source text formatting / precision).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '27.304.02'
$ws.Range("E2").Value = '  +1.97%  '

$ws.Range("D3").Value = '1.660.42'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").Value = '  -0.64%  '

$ws.Range("D5").Value = '220.15'
$ws.Range("E5").Value = '  +1.59%  '

$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("E7").Value = '  -0.68%  '

$ws.Range("D8").Value = '0.255'
$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("D9").Value = '0.0627'
$ws.Range("E9").Value = '  -0.22%  '

$ws.Range("D10").Value = '19.81'
$ws.Range("E10").Value = '  +3.16%  '

$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("D12").Value = '1.889.69'
$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("D13").Value = '1.654.87'
$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("D14").Value = '4.20'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").Value = '0.532'
$ws.Range("E15").Value = '  +0.59%  '

$ws.Range("D16").Value = '66.58'
$ws.Range("E16").Value = '  +2.97%  '

$ws.Range("D17").Value = '27.264.67'
$ws.Range("E17").Value = '  +1.66%  '

$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").Value = '221.60'
$ws.Range("E19").Value = '  +3.45%  '

$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '4.45'
$ws.Range("E21").Value = '  +1.16%  '

$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '6.72'
$ws.Range("E22").Value = '  +7.39%  '

$ws.Range("D23").Value = '2.44'
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").Value = '9.30'
$ws.Range("E24").Value = '  -0.52%  '

$ws.Range("D25").Value = '147.44'

$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("D27").Value = '7.44'
$ws.Range("E27").Value = '  +3.38%  '

$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("D29").Value = '15.98'
$ws.Range("E29").Value = '  +1.92%  '

$ws.Range("E30").Value = '  +1.12%  '

$ws.Range("E31").Value = '  +1.14%  '

$ws.Range("E32").Value = '  +0.57%  '

$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("E34").Value = '  +2.44%  '

$ws.Range("D35").Value = '1.272.14'
$ws.Range("E35").Value = '  -1.17%  '

$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("E37").Value = '  -1.20%  '

$ws.Range("D38").Value = '0.539'
$ws.Range("E38").Value = '  +0.79%  '

$ws.Range("D39").Value = '0.832'
$ws.Range("E39").Value = '  +1.24%  '

$ws.Range("E40").Value = '  -0.50%  '

$ws.Range("D41").Value = '0.810'
$ws.Range("E41").Value = '  +0.47%  '

$ws.Range("D42").Value = '5.38'
$ws.Range("E42").Value = '  +0.92%  '

$ws.Range("D43").Value = '1.801.54'
$ws.Range("E43").Value = '  +0.93%  '

$ws.Range("E44").Value = '  -5.23%  '

$ws.Range("D45").Value = '61.98'
$ws.Range("E45").Value = '  +0.93%  '

$ws.Range("D46").Value = '92.74'
$ws.Range("E46").Value = '  +0.97%  '

$ws.Range("D47").Value = '1.61'
$ws.Range("E47").Value = '  +0.37%  '

$ws.Range("D48").Value = '0.0518'
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").Value = '0.0980'
$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("D50").Value = '7.67'
$ws.Range("E50").Value = '  +0.27%  '

$ws.Range("E51").Value = '  +0.12%  '
